$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row above row 35, shifting existing rows (35-96) down by one.
$ws.Rows.Item(35).Insert()

# Populate the new row's September_Details / September_Date cells.
$ws.Range("R35").Value = "corporate internet share"
$ws.Range("S35").Value = "2024-09-09 11:03:09"
